$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$data = @(
    @("CANADIAN DOLLAR - CHICAGO MERCANTILE EXCHANGE", "11/21/2023", 19407, 84847),
    @("SWISS FRANC - CHICAGO MERCANTILE EXCHANGE", "11/21/2023", 4727, 23721),
    @("BRITISH POUND STERLING - CHICAGO MERCANTILE EXCHANGE", "11/21/2023", 43300, 69398),
    @("JAPANESE YEN - CHICAGO MERCANTILE EXCHANGE", "11/21/2023", 39236, 144690),
    @("U.S. DOLLAR INDEX - ICE FUTURES U.S.", "11/21/2023", 28543, 7846),
    @("EURO FX - CHICAGO MERCANTILE EXCHANGE", "11/21/2023", 231095, 101441),
    @("NEW ZEALAND DOLLAR - CHICAGO MERCANTILE EXCHANGE", "11/21/2023", 7958, 24812),
    @("AUSTRALIAN DOLLAR - CHICAGO MERCANTILE EXCHANGE", "11/21/2023", 29565, 107535),
    @("CANADIAN DOLLAR - CHICAGO MERCANTILE EXCHANGE", "11/28/2023", 18991, 82233),
    @("SWISS FRANC - CHICAGO MERCANTILE EXCHANGE", "11/28/2023", 4300, 24589),
    @("BRITISH POUND STERLING - CHICAGO MERCANTILE EXCHANGE", "11/28/2023", 61296, 69191),
    @("JAPANESE YEN - CHICAGO MERCANTILE EXCHANGE", "11/28/2023", 30461, 139698),
    @("U.S. DOLLAR INDEX - ICE FUTURES U.S.", "11/28/2023", 28798, 9711),
    @("EURO FX - CHICAGO MERCANTILE EXCHANGE", "11/28/2023", 233454, 90289),
    @("NEW ZEALAND DOLLAR - CHICAGO MERCANTILE EXCHANGE", "11/28/2023", 10104, 29713),
    @("AUSTRALIAN DOLLAR - CHICAGO MERCANTILE EXCHANGE", "11/28/2023", 29203, 100422)
)

$startRow = 1450
for ($i = 0; $i -lt $data.Length; $i++) {
    $row = $startRow + $i
    $rec = $data[$i]
    $ws.Cells.Item($row, 1).Value = $rec[0]

    # Column B holds date-like text (e.g. "11/21/2023") that must stay a
    # literal string, not get auto-converted to a date serial number.
    # Temporarily force Text format, assign, then clear the formatting
    # again so the cell keeps the default (unstyled) look.
    $dateCell = $ws.Cells.Item($row, 2)
    $dateCell.NumberFormat = "@"
    $dateCell.Value = $rec[1]
    $dateCell.ClearFormats()

    $ws.Cells.Item($row, 3).Value = $rec[2]
    $ws.Cells.Item($row, 4).Value = $rec[3]
}
